# Updated cryptos list on Tue Mar  7 05:43:15 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) columns for the
# crypto-ranking table on the active sheet. Values are plain display text
# (e.g. "1.152", "22.472.02") rather than real numbers, so every write to
# column D forces a Text number format first and clears it again right
# after so the cell keeps its original (default) style - otherwise Excel
# would happily parse strings like "1.152" as the number 1.152.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Force text storage so numeric-looking strings ("1.152", "0.9999", ...)
    # aren't silently reinterpreted as numbers, then drop the temporary
    # format again so the cell's style index is left untouched.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$updates = @(
    @{ Row = 2;  Price = "22.472.02";  Volume = "  +0.38%  " },
    @{ Row = 3;  Price = "1.575.21";   Volume = "  +0.94%  " },
    @{ Row = 5;  Price = $null;        Volume = "  -0.10%  " },
    @{ Row = 6;  Price = "288.12";     Volume = "  +0.80%  " },
    @{ Row = 7;  Price = "0.3698";     Volume = "  +1.57%  " },
    @{ Row = 8;  Price = "47.73";      Volume = "  -1.98%  " },
    @{ Row = 9;  Price = $null;        Volume = "  -0.09%  " },
    @{ Row = 10; Price = "1.152";      Volume = "  +2.49%  " },
    @{ Row = 11; Price = "0.07575";    Volume = "  +2.67%  " },
    @{ Row = 12; Price = $null;        Volume = "  -0.06%  " },
    @{ Row = 13; Price = $null;        Volume = "  +0.54%  " },
    @{ Row = 14; Price = "5.955";      Volume = "  +0.81%  " },
    @{ Row = 15; Price = $null;        Volume = "  +1.44%  " },
    @{ Row = 16; Price = "1.566.16";   Volume = "  +0.32%  " },
    @{ Row = 17; Price = "0.00001123"; Volume = "  +2.12%  " },
    @{ Row = 18; Price = "88.34";      Volume = "  -0.35%  " },
    @{ Row = 19; Price = "0.06732";    Volume = "  +0.20%  " },
    @{ Row = 20; Price = "0.9999";     Volume = "  -0.13%  " },
    @{ Row = 21; Price = "6.396";      Volume = "  +1.34%  " },
    @{ Row = 22; Price = "16.55";      Volume = "  +3.59%  " },
    @{ Row = 23; Price = "12.05";      Volume = "  +1.25%  " },
    @{ Row = 24; Price = "22.468.90";  Volume = "  +0.37%  " },
    @{ Row = 25; Price = "2.388";      Volume = "  +0.01%  " },
    @{ Row = 26; Price = "2.642";      Volume = "  +3.42%  " },
    @{ Row = 27; Price = "151.26";     Volume = "  +1.40%  " },
    @{ Row = 28; Price = "19.70";      Volume = "  +1.83%  " },
    @{ Row = 29; Price = "4.992";      Volume = "  -0.39%  " },
    @{ Row = 30; Price = "125.60";     Volume = "  +2.35%  " },
    @{ Row = 31; Price = "1.747.47";   Volume = "  +0.69%  " },
    @{ Row = 32; Price = "1.096";      Volume = "  +4.08%  " },
    @{ Row = 33; Price = "6.117";      Volume = "  +0.20%  " },
    @{ Row = 34; Price = "1.986";      Volume = "  -0.12%  " },
    @{ Row = 35; Price = "9.877";      Volume = "  +3.34%  " },
    @{ Row = 36; Price = "0.08377";    Volume = "  +1.92%  " },
    @{ Row = 37; Price = "0.02463";    Volume = "  +4.13%  " },
    @{ Row = 38; Price = "0.2241";     Volume = "  +1.52%  " },
    @{ Row = 39; Price = $null;        Volume = "  +0.93%  " },
    @{ Row = 40; Price = "1.299";      Volume = "  +0.57%  " },
    @{ Row = 41; Price = "5.364";      Volume = "  +1.19%  " },
    @{ Row = 42; Price = $null;        Volume = "  +3.52%  " },
    @{ Row = 43; Price = "0.6285";     Volume = "  +4.09%  " },
    @{ Row = 44; Price = "14.06";      Volume = "  +3.88%  " },
    @{ Row = 45; Price = $null;        Volume = "  -0.10%  " },
    @{ Row = 46; Price = "0.6121";     Volume = "  +7.10%  " },
    @{ Row = 47; Price = $null;        Volume = "  +0.58%  " },
    @{ Row = 48; Price = "2.054";      Volume = "  +2.93%  " },
    @{ Row = 49; Price = "125.33";     Volume = "  +0.80%  " },
    @{ Row = 50; Price = "1.211";      Volume = "  +0.33%  " },
    @{ Row = 51; Price = "0.07220";    Volume = "  +0.05%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        Set-TextValue -range $ws.Range("D$($u.Row)") -text $u.Price
    }
    $ws.Range("E$($u.Row)").Value = $u.Volume
}
